$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the StatQuery formatting (font size 14, wrap text) from C2 onto C3 and C4
$ws.Range("C2").Copy() | Out-Null
$ws.Range("C3:C4").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# C3 and C4 StatQuery cells now use the same full Cypher stat query as C2 (the old short
# CALL{...} direct-Cypher stat query is retired in favor of the parameterized query)
$statQuery = @'
WITH {
    phs_accession: "phs002250",
    subject_ids: [],
    experimental_strategies: [],
    genders: [],
    sample_tumor_statuses: [],
    file_types: [],
    library_strategies: [],
    library_sources: [],
    library_selections: [],
    library_layouts: [],
    platforms: [],
    instrument_models: [],
    reference_genome_assemblies: [],
    primary_diagnoses: [],
    num_study_samples_min: 0,
    num_study_samples_max: 0,
    num_study_participants_max: 0,
    num_study_participants_min: 0
} AS inputs, "Not specified in data" AS na
CALL{
    WITH inputs, na
    MATCH (s:study {phs_accession: inputs.phs_accession})
    OPTIONAL MATCH (s)<--(p:participant)
    WITH inputs, na, s, count(distinct p) AS num_p
    WHERE
        (inputs.num_study_participants_min = 0 OR num_p >= inputs.num_study_participants_min) AND
        (inputs.num_study_participants_max = 0 OR num_p <= inputs.num_study_participants_max)
    OPTIONAL MATCH (s)<--(:participant)<--(samp:sample)
    WITH inputs, na, s, count(distinct samp) AS num_samp
    WHERE
        (inputs.num_study_samples_min = 0 OR num_samp >= inputs.num_study_samples_min) AND
        (inputs.num_study_participants_max = 0 OR num_samp <= inputs.num_study_participants_max)
    MATCH (s)<--(p:participant)
    WITH inputs, na, p
    WHERE
        (size(inputs.subject_ids) = 0 OR p.participant_id IN inputs.subject_ids) AND
        (size(inputs.genders) = 0 OR p.gender IN inputs.genders)
    OPTIONAL MATCH (p)<--(samp:sample)
    WITH inputs, na, p, samp
    WHERE
        (size(inputs.sample_tumor_statuses) = 0 OR (samp is not null AND coalesce(samp.sample_tumor_status, na) IN inputs.sample_tumor_statuses))
    OPTIONAL MATCH (samp)<--(f:file)
    UNWIND coalesce(apoc.text.split(f.experimental_strategy_and_data_subtypes,"[;,]\\s{0,1}"), na) AS experimental_strategies
    WITH inputs, na, p, f, experimental_strategies
    WHERE
        (size(inputs.experimental_strategies) = 0 OR (f is not null AND experimental_strategies IN inputs.experimental_strategies)) AND
        (size(inputs.file_types) = 0 OR (f is not null AND coalesce(f.file_type, na) IN inputs.file_types))
    OPTIONAL MATCH (f)<--(g:genomic_info)
    WITH inputs, na, p, g
    WHERE
        (size(inputs.library_strategies) = 0 OR (g is not null AND coalesce(g.library_strategy , na) IN inputs.library_strategies)) AND
        (size(inputs.library_sources) = 0 OR (g is not null AND coalesce(g.library_source , na) IN inputs.library_sources)) AND
        (size(inputs.library_selections) = 0 OR (g is not null AND coalesce(g.library_selection , na) IN inputs.library_selections)) AND
        (size(inputs.library_layouts) = 0 OR (g is not null AND coalesce(g.library_layout , na) IN inputs.library_layouts)) AND
        (size(inputs.platforms) = 0 OR (g is not null AND coalesce(g.platform , na) IN inputs.platforms)) AND
        (size(inputs.instrument_models) = 0 OR (g is not null AND coalesce(g.instrument_model , na) IN inputs.instrument_models)) AND
        (size(inputs.reference_genome_assemblies) = 0 OR (g is not null AND coalesce(g.reference_genome_assembly , na) IN inputs.reference_genome_assemblies))
    OPTIONAL MATCH (p)<--(diag:diagnosis)
    WITH inputs, na, p, diag
    WHERE
        (size(inputs.primary_diagnoses) = 0 OR (diag is not null AND coalesce(diag.primary_diagnosis, na) IN inputs.primary_diagnoses))
    RETURN
        count(distinct p) AS num_participants
    }
WITH inputs, na, num_participants
CALL {
    WITH inputs, na
    MATCH (s:study {phs_accession: inputs.phs_accession})
    OPTIONAL MATCH (s)<--(p:participant)
    WITH inputs, na, s, count(distinct p) AS num_p
    WHERE
        (inputs.num_study_participants_min = 0 OR num_p >= inputs.num_study_participants_min) AND
        (inputs.num_study_participants_max = 0 OR num_p <= inputs.num_study_participants_max)
    OPTIONAL MATCH (s)<--(:participant)<--(samp:sample)
    WITH inputs, na, s, count(distinct samp) AS num_samp
    WHERE
        (inputs.num_study_samples_min = 0 OR num_samp >= inputs.num_study_samples_min) AND
        (inputs.num_study_participants_max = 0 OR num_samp <= inputs.num_study_participants_max)
    MATCH (s)<--(:participant)<--(samp:sample)
    WITH inputs, na, samp
    WHERE
        (size(inputs.sample_tumor_statuses) = 0 OR coalesce(samp.sample_tumor_status, na) IN inputs.sample_tumor_statuses)
    OPTIONAL MATCH (samp)-->(p:participant)
    WITH inputs, na, samp, p
    WHERE
        (size(inputs.subject_ids) = 0 OR (p is not null AND p.participant_id IN inputs.subject_ids)) AND
        (size(inputs.genders) = 0 OR (p is not null AND p.gender IN inputs.genders))
    OPTIONAL MATCH (samp)<--(f:file)
    UNWIND coalesce(apoc.text.split(f.experimental_strategy_and_data_subtypes,"[;,]\\s{0,1}"), na) AS experimental_strategies
    WITH inputs, na, samp, f, experimental_strategies
    WHERE
        (size(inputs.experimental_strategies) = 0 OR (f is not null AND experimental_strategies IN inputs.experimental_strategies)) AND
        (size(inputs.file_types) = 0 OR (f is not null AND coalesce(f.file_type, na) IN inputs.file_types))
    OPTIONAL MATCH (f)<--(g:genomic_info)
    WITH inputs, na, samp, g
    WHERE
        (size(inputs.library_strategies) = 0 OR (g is not null AND coalesce(g.library_strategy , na) IN inputs.library_strategies)) AND
        (size(inputs.library_sources) = 0 OR (g is not null AND coalesce(g.library_source , na) IN inputs.library_sources)) AND
        (size(inputs.library_selections) = 0 OR (g is not null AND coalesce(g.library_selection , na) IN inputs.library_selections)) AND
        (size(inputs.library_layouts) = 0 OR (g is not null AND coalesce(g.library_layout , na) IN inputs.library_layouts)) AND
        (size(inputs.platforms) = 0 OR (g is not null AND coalesce(g.platform , na) IN inputs.platforms)) AND
        (size(inputs.instrument_models) = 0 OR (g is not null AND coalesce(g.instrument_model , na) IN inputs.instrument_models)) AND
        (size(inputs.reference_genome_assemblies) = 0 OR (g is not null AND coalesce(g.reference_genome_assembly , na) IN inputs.reference_genome_assemblies))
    OPTIONAL MATCH (samp)-->(:participant)<--(diag:diagnosis)
    WITH inputs, na, samp, diag
    WHERE
        (size(inputs.primary_diagnoses) = 0 OR (diag is not null AND coalesce(diag.primary_diagnosis, na) IN inputs.primary_diagnoses))
    RETURN
        count(distinct samp) AS num_samples
}
WITH inputs, na, num_participants, num_samples
CALL {
    WITH inputs, na
    MATCH (s:study {phs_accession: inputs.phs_accession})
    OPTIONAL MATCH (s)<--(p:participant)
    WITH inputs, na, s, count(distinct p) AS num_p
    WHERE
        (inputs.num_study_participants_min = 0 OR num_p >= inputs.num_study_participants_min) AND
        (inputs.num_study_participants_max = 0 OR num_p <= inputs.num_study_participants_max)
    OPTIONAL MATCH (s)<--(:participant)<--(samp:sample)
    WITH inputs, na, s, count(distinct samp) AS num_samp
    WHERE
        (inputs.num_study_samples_min = 0 OR num_samp >= inputs.num_study_samples_min) AND
        (inputs.num_study_participants_max = 0 OR num_samp <= inputs.num_study_participants_max)
    MATCH (s)<--(f:file)
    UNWIND coalesce(apoc.text.split(f.experimental_strategy_and_data_subtypes,"[;,]\\s{0,1}"), na) AS experimental_strategies
    WITH inputs, na, f, experimental_strategies
    WHERE
        (size(inputs.experimental_strategies) = 0 OR experimental_strategies IN inputs.experimental_strategies) AND
        (size(inputs.file_types) = 0 OR coalesce(f.file_type, na) IN inputs.file_types)
    OPTIONAL MATCH (f)-->(samp:sample)
    WITH inputs, na, f, samp
    WHERE
        (size(inputs.sample_tumor_statuses) = 0 OR (samp is not null AND coalesce(samp.sample_tumor_status, na) IN inputs.sample_tumor_statuses))
    OPTIONAL MATCH (samp)-->(p:participant)
    WITH inputs, na, f, p
    WHERE
        (size(inputs.subject_ids) = 0 OR (p is not null AND p.participant_id IN inputs.subject_ids)) AND
        (size(inputs.genders) = 0 OR (p is not null AND p.gender IN inputs.genders))
    OPTIONAL MATCH (f)<--(g:genomic_info)
    WITH inputs, na, f, g
    WHERE
        (size(inputs.library_strategies) = 0 OR (g is not null AND coalesce(g.library_strategy , na) IN inputs.library_strategies)) AND
        (size(inputs.library_sources) = 0 OR (g is not null AND coalesce(g.library_source , na) IN inputs.library_sources)) AND
        (size(inputs.library_selections) = 0 OR (g is not null AND coalesce(g.library_selection , na) IN inputs.library_selections)) AND
        (size(inputs.library_layouts) = 0 OR (g is not null AND coalesce(g.library_layout , na) IN inputs.library_layouts)) AND
        (size(inputs.platforms) = 0 OR (g is not null AND coalesce(g.platform , na) IN inputs.platforms)) AND
        (size(inputs.instrument_models) = 0 OR (g is not null AND coalesce(g.instrument_model , na) IN inputs.instrument_models)) AND
        (size(inputs.reference_genome_assemblies) = 0 OR (g is not null AND coalesce(g.reference_genome_assembly , na) IN inputs.reference_genome_assemblies))
    OPTIONAL MATCH (f)-->(:sample)-->(:participant)<--(diag:diagnosis)
    WITH inputs, na, f, diag
    WHERE
        (size(inputs.primary_diagnoses) = 0 OR coalesce(diag.primary_diagnosis, na) IN inputs.primary_diagnoses)
    RETURN
        count(distinct f) AS num_files
}
RETURN
    1 AS Studies,
    num_participants AS Participants,
    num_samples AS Samples,
    num_files AS `Files`
'@
$ws.Range("C3").Value = $statQuery
$ws.Range("C4").Value = $statQuery

# FilesTab query: clear the RNA-Seq experimental_strategies filter back to an empty list
$filesQuery = @'
WITH {
    phs_accession: "phs002250",
    subject_ids: [],
    experimental_strategies: [],
    genders: [],
    sample_tumor_statuses: [],
    file_types: [],
    library_strategies: [],
    library_sources: [],
    library_selections: [],
    library_layouts: [],
    platforms: [],
    instrument_models: [],
    reference_genome_assemblies: [],
    primary_diagnoses: [],
    num_study_samples_min: 0,
    num_study_samples_max: 0,
    num_study_participants_max: 0,
    num_study_participants_min: 0
} AS inputs, "Not specified in data" AS na
MATCH (s:study {phs_accession: inputs.phs_accession})
OPTIONAL MATCH (s)<--(p:participant)
WITH inputs, na, s, count(distinct p) AS num_p
WHERE
    (inputs.num_study_participants_min = 0 OR num_p >= inputs.num_study_participants_min) AND
    (inputs.num_study_participants_max = 0 OR num_p <= inputs.num_study_participants_max)
OPTIONAL MATCH (s)<--(:participant)<--(samp:sample)
WITH inputs, na, s, count(distinct samp) AS num_samp
WHERE
    (inputs.num_study_samples_min = 0 OR num_samp >= inputs.num_study_samples_min) AND
    (inputs.num_study_participants_max = 0 OR num_samp <= inputs.num_study_participants_max)
MATCH (s)<--(f:file)
WITH inputs, na, f, {
file_name: coalesce(f.file_name, ""),
file_type: coalesce(f.file_type, ""),
study_name: coalesce(s.study_name, ""),
accession: coalesce(s.phs_accession, "")
} AS data,
apoc.coll.toSet(apoc.text.split(f.experimental_strategy_and_data_subtypes,"[;,]\\s{0,1}")) AS experimental_strategies
WHERE
(size(inputs.experimental_strategies) = 0 OR size(apoc.coll.intersection(inputs.experimental_strategies, experimental_strategies)) > 0) AND
    (size(inputs.file_types) = 0 OR f.file_type IN inputs.file_types)
OPTIONAL MATCH (f)-->(samp)
WITH inputs, na, f, apoc.map.merge(data, {
sample_id: coalesce(apoc.text.join(apoc.coll.sort(collect(distinct samp.sample_id)), ", "), "")
}) AS data,
COLLECT(DISTINCT samp.sample_tumor_status) AS sample_tumor_status
WHERE
(size(inputs.sample_tumor_statuses) = 0 OR size(apoc.coll.intersection(inputs.sample_tumor_statuses, sample_tumor_status)) > 0)
OPTIONAL MATCH (f)-->(:sample)-->(p:participant)
WITH inputs, na, f,
apoc.map.merge(data, {
participant_id: coalesce(apoc.text.join(collect(distinct p.participant_id), ", "), "")
}) AS data,
COLLECT(DISTINCT p.participant_id) AS participant_id,
COLLECT(DISTINCT p.gender) AS gender
WHERE
    (size(inputs.subject_ids) = 0 OR size(apoc.coll.intersection(inputs.participant_id, participant_id)) > 0) AND
    (size(inputs.genders) = 0 OR size(apoc.coll.intersection(inputs.gender, gender)) > 0)
OPTIONAL MATCH (f)<--(g:genomic_info)
WITH inputs, na, f, data,
COLLECT(DISTINCT g.library_strategy) AS library_strategy,
COLLECT(DISTINCT g.library_source) AS library_source,
COLLECT(DISTINCT g.library_selection) AS library_selection,
COLLECT(DISTINCT g.library_layout) AS library_layout,
COLLECT(DISTINCT g.platform) AS platform,
COLLECT(DISTINCT g.instrument_model) AS instrument_models,
COLLECT(DISTINCT g.reference_genome_assembly) AS reference_genome_assembly
WHERE
    (size(inputs.library_strategies) = 0 OR size(apoc.coll.intersection(inputs.library_strategies, library_strategy)) > 0) AND
    (size(inputs.library_sources) = 0 OR size(apoc.coll.intersection(inputs.library_sources, library_source)) > 0) AND
    (size(inputs.library_selections) = 0 OR size(apoc.coll.intersection(inputs.library_selections, library_selection)) > 0) AND
    (size(inputs.library_layouts) = 0 OR size(apoc.coll.intersection(inputs.library_layouts, library_layout)) > 0) AND
    (size(inputs.platforms) = 0 OR size(apoc.coll.intersection(inputs.platforms, platform)) > 0) AND
    (size(inputs.instrument_models) = 0 OR size(apoc.coll.intersection(inputs.instrument_models, instrument_models)) > 0) AND
    (size(inputs.reference_genome_assemblies) = 0 OR size(apoc.coll.intersection(inputs.reference_genome_assemblies, reference_genome_assembly)) > 0)
OPTIONAL MATCH (f)-->(:sample)-->(:participant)<--(diag:diagnosis)
WITH inputs, na, f, data,
COLLECT(DISTINCT diag.primary_diagnosis) AS primary_diagnosis
WHERE
    (size(inputs.primary_diagnoses) = 0 OR size(apoc.coll.intersection(inputs.primary_diagnoses, primary_diagnosis)) > 0)
RETURN
data.file_name AS `File Name`,
data.study_name AS `Study Name`,
data.accession AS `Accession`,
data.participant_id AS `Participant Id`,
data.sample_id AS `Sample Id`,
data.file_type AS `File Type`
ORDER BY `File Name`
LIMIT 100
'@
$ws.Range("B4").Value = $filesQuery

# Re-assert the (already maxed-out) row heights; re-entering the long wrapped
# text above can nudge Excel's auto row-height math, so pin them back to the
# sheet's real max height.
$ws.Rows.Item(2).RowHeight = 409.5
$ws.Rows.Item(3).RowHeight = 409.5
$ws.Rows.Item(4).RowHeight = 409.5

# Restore the selection Excel leaves on this tab after the edit
$ws.Range("C4").Select()

